# "changed from qos 2 to qos 1 and excel sheet updated"
#
# The throughput results for the MQTT QoS1 run (row 4) and the MQTT QoS2
# run (row 5) were previously blank placeholders; fill them in with the
# measured values now that the run has been completed/updated. Also widen
# the data columns to comfortably show the new (much longer) numbers, and
# leave the active cell where the user ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 - MQTT QoS1
$ws.Range("B4").Value = 83539.399999999994
$ws.Range("C4").Value = 17885.09
$ws.Range("D4").Value = 828635.17299999995
$ws.Range("E4").Value = 141518.14199999999
$ws.Range("F4").Value = 844990176.227
$ws.Range("G4").Value = 170006467.08500001
$ws.Range("H4").Value = 7041138301.8100004
$ws.Range("I4").Value = 3214330801.8800001

# Row 5 - MQTT QoS2
$ws.Range("B5").Value = 1668.6681000000001
$ws.Range("C5").Value = 682.29
$ws.Range("D5").Value = 17726.009999999998
$ws.Range("E5").Value = 7938.2
$ws.Range("F5").Value = 46784083.530000001
$ws.Range("G5").Value = 16843543.767000001
$ws.Range("H5").Value = 511400284.29500002
$ws.Range("I5").Value = 222389616.449

# Widen the columns that now hold the newly-entered figures (and the two
# narrow trailing columns) so the values aren't clipped/##### on screen.
$ws.Columns.Item(1).ColumnWidth = 10.251822916666665
$ws.Columns.Item(2).ColumnWidth = 13.918489583333333
$ws.Columns.Item(3).ColumnWidth = 13.251822916666665
$ws.Columns.Item(4).ColumnWidth = 12.918489583333333
$ws.Columns.Item(5).ColumnWidth = 13.085156249999999
$ws.Columns.Item(6).ColumnWidth = 14.418489583333333
$ws.Columns.Item(7).ColumnWidth = 13.585156249999999
$ws.Columns.Item(8).ColumnWidth = 20.25182291666667
$ws.Columns.Item(9).ColumnWidth = 21.418489583333333
$ws.Columns.Item(14).ColumnWidth = 8.251822916666665
$ws.Columns.Item(15).ColumnWidth = 8.251822916666665
$ws.Columns.Item(16).ColumnWidth = 8.251822916666665
$ws.Columns.Item(17).ColumnWidth = 8.251822916666665

# Move the active selection to where the user was last working.
$ws.Range("H4").Select()
